# Update countries & provincias Spain
# Applies the 28-Abril-2020 00:22 data refresh to the "Pais" sheet:
#   - refreshed case counts for several countries
#   - re-ranking swaps for Tunez/Bolivia, Gabon/Ruanda/Congo/Islas Feroe,
#     and Bermudas/Cabo Verde (values cascade down one row, with the
#     newly top-ranked country getting the fresh totals)
#   - updated "Datos actualizados a ..." timestamp

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp banner ---
$ws.Range("A1").Value = "Datos actualizados a 28 de Abril de 2020 a las 00:22"

# --- Estados Unidos (row 4) ---
$ws.Range("B4").Value = 1005808
$ws.Range("C4").Value = 18648
$ws.Range("D4").Value = 137693
$ws.Range("E4").Value = 811554
$ws.Range("F4").Value = 14175
$ws.Range("G4").Value = 1148
$ws.Range("H4").Value = 56561

# --- Canada (row 15) ---
$ws.Range("B15").Value = 48458
$ws.Range("C15").Value = 1563
$ws.Range("D15").Value = 18215
$ws.Range("E15").Value = 27539
$ws.Range("F15").Value = 557
$ws.Range("G15").Value = 144
$ws.Range("H15").Value = 2704

# --- Tunez overtakes Bolivia (rows 90-91) ---
$ws.Range("A90").Value = "Tunez"
$ws.Range("B90").Value = 967
$ws.Range("C90").Value = 18
$ws.Range("D90").Value = 279
$ws.Range("E90").Value = 649
$ws.Range("F90").Value = 18
$ws.Range("G90").Value = 1
$ws.Range("H90").Value = 39

$ws.Range("A91").Value = "Bolivia"
$ws.Range("B91").Value = 950
$ws.Range("C91").Value = 84
$ws.Range("D91").Value = 80
$ws.Range("E91").Value = 820
$ws.Range("F91").Value = 3
$ws.Range("G91").Value = 4
$ws.Range("H91").Value = 50

# --- Gabon overtakes Ruanda, Congo, Islas Feroe (rows 131-134) ---
$ws.Range("A131").Value = "Gabon"
$ws.Range("B131").Value = 211
$ws.Range("C131").Value = 35
$ws.Range("D131").Value = 43
$ws.Range("E131").Value = 165
$ws.Range("F131").Value = 1
$ws.Range("G131").Value = 0
$ws.Range("H131").Value = 3

$ws.Range("A132").Value = "Ruanda"
$ws.Range("B132").Value = 207
$ws.Range("C132").Value = 16
$ws.Range("D132").Value = 93
$ws.Range("E132").Value = 114
$ws.Range("F132").Value = 0
$ws.Range("G132").Value = 0
$ws.Range("H132").Value = 0

$ws.Range("A133").Value = "Congo"
$ws.Range("B133").Value = 200
$ws.Range("C133").Value = 0
$ws.Range("D133").Value = 19
$ws.Range("E133").Value = 175
$ws.Range("F133").Value = 0
$ws.Range("G133").Value = 0
$ws.Range("H133").Value = 6

$ws.Range("A134").Value = "Islas Feroe"
$ws.Range("B134").Value = 187
$ws.Range("C134").Value = 0
$ws.Range("D134").Value = 178
$ws.Range("E134").Value = 9
$ws.Range("F134").Value = 0
$ws.Range("G134").Value = 0
$ws.Range("H134").Value = 0

# --- Bermudas overtakes Cabo Verde (rows 146-147) ---
$ws.Range("A146").Value = "Bermudas"
$ws.Range("B146").Value = 110
$ws.Range("C146").Value = 1
$ws.Range("D146").Value = 44
$ws.Range("E146").Value = 60
$ws.Range("F146").Value = 10
$ws.Range("G146").Value = 1
$ws.Range("H146").Value = 6

$ws.Range("A147").Value = "Cabo Verde"
$ws.Range("B147").Value = 109
$ws.Range("C147").Value = 3
$ws.Range("D147").Value = 1
$ws.Range("E147").Value = 107
$ws.Range("F147").Value = 0
$ws.Range("G147").Value = 0
$ws.Range("H147").Value = 1
